$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 5.552368248186704
$ws.Range("C2").Value = 13156.8
$ws.Range("D2").Value = 0.03811254470716262
$ws.Range("E2").Value = 92.2
$ws.Range("F2").Value = 208.8
$ws.Range("G2").Value = "Inuwâ Dân X"
$ws.Range("H2").Value = "DUO"
$ws.Range("I2").Value = 0.08970024567781631
$ws.Range("J2").Value = 10.2
$ws.Range("K2").Value = 0.0041437734423764
$ws.Range("L2").Value = 5.4
$ws.Range("M2").Value = 2.8
$ws.Range("N2").Value = 4.8
$ws.Range("O2").Value = 15.2
$ws.Range("P2").Value = 13156.8
$ws.Range("Q2").Value = 10.2
$ws.Range("B3").Value = 8.640537996803857
$ws.Range("C3").Value = 19508.6
$ws.Range("D3").Value = 0.06836728990692152
$ws.Range("E3").Value = 110.8
$ws.Range("F3").Value = 342.8
$ws.Range("G3").Value = "Takizawà"
$ws.Range("H3").Value = "DUO_CARRY"
$ws.Range("I3").Value = 0.1609056249067745
$ws.Range("J3").Value = 19.2
$ws.Range("K3").Value = 0.008901090469941059
$ws.Range("L3").Value = 3.2
$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 7.6
$ws.Range("O3").Value = 14
$ws.Range("P3").Value = 19508.6
$ws.Range("Q3").Value = 19.2
$ws.Range("B4").Value = 4.958948577935138
$ws.Range("C4").Value = 7422.2
$ws.Range("D4").Value = 0.05297472115560431
$ws.Range("E4").Value = 75.8
$ws.Range("F4").Value = 343.6
$ws.Range("G4").Value = "1against9 L9"
$ws.Range("H4").Value = "DUO_CARRY"
$ws.Range("I4").Value = 0.2438613715321066
$ws.Range("J4").Value = 4.8
$ws.Range("K4").Value = 0.003429869180006317
$ws.Range("L4").Value = 1.8
$ws.Range("M4").Value = 1.2
$ws.Range("N4").Value = 2.2
$ws.Range("O4").Value = 1.4
$ws.Range("P4").Value = 7422.2
$ws.Range("Q4").Value = 4.8
$ws.Range("B5").Value = 7.578175204713018
$ws.Range("C5").Value = 13639.6
$ws.Range("D5").Value = 0.07927726354393795
$ws.Range("E5").Value = 145
$ws.Range("F5").Value = 588.4
$ws.Range("G5").Value = "Mr F Scofield"
$ws.Range("H5").Value = "DUO_CARRY"
$ws.Range("I5").Value = 0.3224582848240293
$ws.Range("J5").Value = 12
$ws.Range("K5").Value = 0.006465115763260081
$ws.Range("L5").Value = 6.2
$ws.Range("M5").Value = 3.2
$ws.Range("N5").Value = 6
$ws.Range("O5").Value = 6.2
$ws.Range("P5").Value = 13639.6
$ws.Range("Q5").Value = 12
$ws.Range("B6").Value = 4.718499124373624
$ws.Range("C6").Value = 6698
$ws.Range("D6").Value = 0.03749183615672267
$ws.Range("E6").Value = 52.4
$ws.Range("F6").Value = 260.8
$ws.Range("G6").Value = "atkı"
$ws.Range("H6").Value = "DUO_CARRY"
$ws.Range("I6").Value = 0.19180321236396
$ws.Range("J6").Value = 7.2
$ws.Range("K6").Value = 0.005053780220669405
$ws.Range("L6").Value = 3.8
$ws.Range("M6").Value = 1.4
$ws.Range("N6").Value = 2.6
$ws.Range("O6").Value = 11.4
$ws.Range("P6").Value = 6698
$ws.Range("Q6").Value = 7.2
